$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end ---
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"
$sheetCount = $wb.Worksheets.Count
$wsForecast.Move($null, $wb.Worksheets.Item($sheetCount))

# Re-fetch the sheet by name since the previous handle now points at
# whatever sheet occupies that old position after the Move.
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows ---
$poForecastData = @(
  @(44997.99999999999, 9, -15.7018589171405, 32.77969645411343),
  @(45004.99999999999, 10, -13.20895278813307, 32.77950628249428),
  @(45039.99999999999, 11, -12.64463532040394, 36.45090449043059),
  @(45046.99999999999, 12, -12.51430315899573, 36.02591848595004),
  @(45053.99999999999, 12, -13.12855618025969, 36.75857029154984),
  @(45067.99999999999, 13, -11.2968933855103, 37.89524649695689),
  @(45081.99999999999, 13, -9.497045750314603, 38.52171439705794),
  @(45088.99999999999, 14, -12.13632536973797, 38.26659933845),
  @(45095.99999999999, 14, -9.245194303530727, 39.00761187573919),
  @(45389.99999999999, 28, 3.834707140396134, 52.73140462531593),
  @(45410.99999999999, 29, 2.693378859312622, 52.41384360992265),
  @(45417.99999999999, 29, 4.927016679810935, 53.93467002973514),
  @(45424.99999999999, 30, 5.063163325605875, 52.56069052004861),
  @(45438.99999999999, 30, 5.436355506606902, 54.82337616088536),
  @(45445.99999999999, 31, 6.146155829436024, 55.03793225531241),
  @(45452.99999999999, 31, 7.982690293499312, 54.93291547791478),
  @(45466.99999999999, 32, 4.824308519252039, 55.27692724531126),
  @(45473.99999999999, 32, 6.294311411525282, 56.02159406698306),
  @(45487.99999999999, 33, 9.767917980110674, 57.01724091301145),
  @(45494.99999999999, 33, 7.637004629599104, 55.92005198199481),
  @(45501.99999999999, 33, 7.752180554492951, 56.13985412538223),
  @(45508.99999999999, 34, 9.29835849118405, 58.43850183269878),
  @(45515.99999999999, 34, 10.2941017853416, 59.4613635527009),
  @(45522.99999999999, 34, 9.548410694912056, 57.86794889961192),
  @(45529.99999999999, 35, 12.66670674879166, 58.65602764234782),
  @(45550.99999999999, 36, 11.43614150295956, 61.41884509133769),
  @(45557.99999999999, 36, 10.51769889585408, 59.3623643486406),
  @(45564.99999999999, 36, 12.76458045995765, 63.03166080489345),
  @(45571.99999999999, 37, 12.8970517779563, 58.13808143231788),
  @(45578.99999999999, 37, 13.28948200443731, 60.61846699785233),
  @(45585.99999999999, 37, 11.66896692745604, 60.82231155924421),
  @(45599.99999999999, 38, 14.80820551466769, 61.5428224526169),
  @(45606.99999999999, 38, 16.75181950481067, 64.24939784214072),
  @(45613.99999999999, 39, 14.21566156510182, 62.74438216235945),
  @(45620.99999999999, 39, 15.75852900086581, 62.82570391331637),
  @(45627.99999999999, 39, 14.51220567197227, 64.19609726887926),
  @(45634.99999999999, 40, 16.49119852921687, 64.78115971878776),
  @(45641.99999999999, 40, 16.41785101161013, 62.84088789469779),
  @(45648.99999999999, 41, 16.11736342752426, 64.36744379836371),
  @(45655.99999999999, 41, 17.77182956333502, 64.71131503915504)
)

for ($i = 0; $i -lt $poForecastData.Count; $i++) {
    $r = $i + 2
    $rowVals = $poForecastData[$i]
    $wsForecast.Cells.Item($r, 1).Value = $rowVals[0]
    $wsForecast.Cells.Item($r, 2).Value = $rowVals[1]
    $wsForecast.Cells.Item($r, 3).Value = $rowVals[2]
    $wsForecast.Cells.Item($r, 4).Value = $rowVals[3]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
